$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - habitat_type
$ws.Range("B2").Value = 207.316926095507
$ws.Range("D2").Value = 0.000000000000000000000000000000000000000000100340404728866

# Row 3 - season
$ws.Range("B3").Value = 44.8659925864483
$ws.Range("D3").Value = 0.000000000987978216148026

# Row 4 - day_night (only statistic changes, p.value stays 0)
$ws.Range("B4").Value = 1665.54589547028

# Row 5 - habitat_type:season
$ws.Range("B5").Value = 267.455553278609
$ws.Range("D5").Value = 0.0000000000000000000000000000000000000000000000000309821676830582

# Row 6 - habitat_type:day_night
$ws.Range("B6").Value = 149.330142220592
$ws.Range("D6").Value = 0.0000000000000000000000000775257765445659

# Row 7 - season:day_night
$ws.Range("B7").Value = 143.019684211909
$ws.Range("D7").Value = 0.0000000000000000000000000245347438577342

# Row 8 - habitat_type:season:day_night
$ws.Range("B8").Value = 123.832365414437
$ws.Range("D8").Value = 0.00000000000747450162815215
